$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = "id"
$ws.Range("H12").Value = "ifelse(is.na(cignr_day),`r`n                      ifelse(is.na(cigarnr_day), `r`n                             ifelse(is.na(pipe_day),0,pipe_day*5),`r`n                             ifelse(is.na(pipe_day),cigarnr_day*5,cigarnr_day*5+pipe_day*5)),`r`n                      ifelse(is.na(cigarnr_day),`r`n                             ifelse(is.na(pipe_day),cignr_day*1,cignr_day*1+pipe_day*5),`r`n                             ifelse(is.na(pipe_day),cignr_day*1+cigarnr_day*5, cignr_day*1+ cigarnr_day*5 + pipe_day*5)))"
$ws.Range("H15").Value = "case_when(`r`n  fr_op7 == 2 ~ 3, `r`n  fr_period1 == 3 ~ 1,`r`n  fr_period1 == 1 ~ 0,`r`n  fr_period1 == 2 & age >= 45 & age <= 55 ~ 2,`r`n  fr_period1 == 2 & age < 45 ~ 0,`r`n  fr_period1 == 2 & age > 55 ~ 1,`r`n  TRUE ~ NA_integer_`r`n)"
$ws.Range("I15").Value = "Perimenopausal category defined as `"yes, irregular`" period and age between 45 and 55 according to NAKO-QS-Report NAKO-1041"
$ws.Range("G25").Value = "case_when"
$ws.Range("H25").Value = "case_when(                                                                                                                                                                                                                                                                                                                                                                                                                                `r`n  hyp_i == 1 ~ 1, `r`n  hyp_i == 2 ~ 0, `r`n  is.na(hyp_i) & htn_kora == 1 ~ 1,`r`n  is.na(hyp_i) & htn_kora == 0 ~ 0, `r`n  TRUE ~ NA_integer_)"
$ws.Range("H28").Value = "recode(1= 1; 2 = 0; 3=2; 8=2;)"
$ws.Range("I28").Value = "Category `"No / I don't know`" will be harmonised to `"I don't know`""
$ws.Range("H29").Value = "recode(1= 1; 2 = 0; 3=2; 8=2;)"
$ws.Range("I29").Value = "Category `"No / I don't know`" will be harmonised to `"I don't know`""
$ws.Range("H30").Value = "recode(1= 1; 2 = 0; 3=2; 8=2;)"
$ws.Range("I30").Value = "Category `"No / I don't know`" will be harmonised to `"I don't know`""
$ws.Range("F38").Value = "med_stat"
$ws.Range("F39").Value = "med_nsaid"
$ws.Range("H52").Value = "case_when(`r`n      f1_htn_kora == 1 | f2_htn_kora == 1 ~ 1,`r`n      f1_htn_kora == 0 & f2_htn_kora == 0 ~ 0,`r`n      TRUE ~ NA_integer_`r`n    )"
$ws.Range("F53").Value = "f1_htn_kora;f2_htn_kora;f1_untdat;f2_untdat;gebdat"
$ws.Range("H53").Value = "case_when(`r`n  f1_htn_kora == 1 ~ as.numeric(f1_untdat - gebdat)/365.25,              `r`n  f2_htn_kora == 1 ~ as.numeric(f2_untdat - gebdat)/365.25,  `r`n  TRUE ~ NA_real_`r`n)"
$ws.Range("H55").Value = "(f2_insuff_date-gebdat)/365.25"
$ws.Range("H62").Value = "(tod_dat-gebdat)/365.25"
